# Apply scheduled-runner updates to Sheets (per commit: 'chore: update Sheets via scheduled runner')
$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$updates = @{
    "H17" = 1270.409
    "J17" = 1218.4762
    "L17" = 3655.4286
    "N17" = -3991.4286
    "H42" = 90.42856999999999
    "I42" = 54
    "K42" = 162
    "M42" = 68
    "H64" = 10544.75
    "I64" = 5639.5
    "J64" = 12179.833
    "K64" = 5639.5
    "L64" = 12179.833
    "M64" = -5391.5
    "N64" = -12675.833
    "H67" = 10544.75
    "I67" = 5639.5
    "J67" = 12179.833
    "K67" = 5639.5
    "L67" = 12179.833
    "M67" = -4781.5
    "N67" = -13895.833
    "H92" = 923.9
    "I92" = 953.8
    "K92" = 953.8
    "M92" = 294.2
    "H112" = 1704.2354
    "J112" = 1705.5938
    "L112" = 5116.7814
    "N112" = -7332.7814
    "H130" = 90000
    "J130" = 90000
    "L130" = 90000
    "N130" = -100040
    "H132" = 14471.619
    "J132" = 60576.332
    "L132" = 181728.996
    "N132" = -186788.996
    "H134" = 70000
    "J134" = 70000
    "L134" = 70000
    "N134" = -80140
    "H135" = 2387
    "I135" = 2346.1333
    "J135" = 3000
    "K135" = 21115.1997
    "L135" = 27000
    "M135" = -18580.1997
    "N135" = -32070
    "H138" = 3552.3125
    "I138" = 2618
    "J138" = 3863.75
    "K138" = 7854
    "L138" = 11591.25
    "M138" = -2714
    "N138" = -21871.25
}
foreach ($cell in $updates.Keys) {
    $ws.Range($cell).Value = $updates[$cell]
}

$ws = $wb.Worksheets.Item("ARM")
$updates = @{
    "H31" = 13023.77
    "I31" = 11380
    "J31" = 18503
    "K31" = 11380
    "L31" = 18503
    "M31" = -11086
    "N31" = -19091
    "H32" = 4946.5415
    "I32" = 5060.8843
    "K32" = 5060.8843
    "M32" = -4773.8843
    "H45" = 2056.9412
    "I45" = 1228.3846
    "K45" = 1228.3846
    "M45" = -851.3846000000001
    "H61" = 1465.5333
    "I61" = 1213.0714
    "K61" = 1213.0714
    "M61" = -1001.0714
    "H74" = 1799.9354
    "I74" = 1825.9584
    "J74" = 1710.7142
    "K74" = 1825.9584
    "L74" = 1710.7142
    "M74" = -951.9584
    "N74" = -3458.7142
    "H77" = 1799.9354
    "I77" = 1825.9584
    "J77" = 1710.7142
    "K77" = 9129.791999999999
    "L77" = 8553.571
    "M77" = -4761.791999999999
    "N77" = -17289.571
    "H88" = 902.5263
    "I88" = 1147.6666
    "J88" = 681.9
    "K88" = 1147.6666
    "L88" = 681.9
    "M88" = -741.6666
    "N88" = -1493.9
    "H91" = 902.5263
    "I91" = 1147.6666
    "J91" = 681.9
    "K91" = 1147.6666
    "L91" = 681.9
    "M91" = 256.3334
    "N91" = -3489.9
    "H102" = 4137.2666
    "I102" = 4137.2666
    "K102" = 4137.2666
    "M102" = -2515.2666
    "H136" = 1465.5333
    "I136" = 1213.0714
    "K136" = 3639.2142
    "M136" = -1089.2142
}
foreach ($cell in $updates.Keys) {
    $ws.Range($cell).Value = $updates[$cell]
}

$ws = $wb.Worksheets.Item("BSM")
$updates = @{
    "H94" = 2352.1562
    "J94" = 3175.6667
    "L94" = 3175.6667
    "N94" = -4077.6667
    "H99" = 68973.664
    "I99" = 68973.664
    "K99" = 68973.664
    "M99" = -67475.664
    "H105" = 2391.8
    "I105" = 2003
    "K105" = 2003
    "M105" = -256
    "H107" = 7361.6665
    "I107" = 6377.846
    "J107" = 13756.5
    "K107" = 6377.846
    "L107" = 13756.5
    "M107" = -4457.846
    "N107" = -17596.5
    "H109" = 29899.5
    "J109" = 29899.5
    "L109" = 29899.5
    "N109" = -32673.5
}
foreach ($cell in $updates.Keys) {
    $ws.Range($cell).Value = $updates[$cell]
}

$ws = $wb.Worksheets.Item("CRP")
$updates = @{
    "H31" = 1862.826
    "I31" = 1767.7059
    "J31" = 2132.3333
    "K31" = 1767.7059
    "L31" = 2132.3333
    "M31" = -1472.7059
    "N31" = -2722.3333
    "H34" = 1862.826
    "I34" = 1767.7059
    "J34" = 2132.3333
    "K34" = 1767.7059
    "L34" = 2132.3333
    "M34" = -1565.7059
    "N34" = -2536.3333
    "H58" = 2483.077
    "I58" = 1468.5714
    "J58" = 3666.6667
    "K58" = 1468.5714
    "L58" = 3666.6667
    "M58" = -1265.5714
    "N58" = -4072.6667
    "H99" = 3193.92
    "J99" = 3111.1538
    "L99" = 3111.1538
    "N99" = -6107.1538
    "H105" = 468.5625
    "I105" = 433.13333
    "K105" = 433.13333
    "M105" = 1313.86667
    "H126" = 3193.92
    "J126" = 3111.1538
    "L126" = 9333.4614
    "N126" = -14273.4614
    "H136" = 2483.077
    "I136" = 1468.5714
    "J136" = 3666.6667
    "K136" = 4405.7142
    "L136" = 11000.0001
    "M136" = -1855.7142
    "N136" = -16100.0001
}
foreach ($cell in $updates.Keys) {
    $ws.Range($cell).Value = $updates[$cell]
}

$ws = $wb.Worksheets.Item("CUL")
$updates = @{
    "H50" = 858.5
    "J50" = 1300
    "L50" = 3900
    "N50" = -4862
    "H53" = 858.5
    "J53" = 1300
    "L53" = 3900
    "N53" = -4862
    "H109" = 820.5
    "I109" = 820.5
    "K109" = 2461.5
    "M109" = -1421.5
    "H118" = 4888.5557
    "I118" = 4332.8335
    "K118" = 12998.5005
    "M118" = -11755.5005
    "H121" = 1815.6666
    "J121" = 2223.625
    "L121" = 6670.875
    "N121" = -9290.875
    "H131" = 3054.0557
    "I131" = 1557.5454
    "K131" = 4672.6362
    "M131" = 367.3638000000001
}
foreach ($cell in $updates.Keys) {
    $ws.Range($cell).Value = $updates[$cell]
}

$ws = $wb.Worksheets.Item("GSM")
$updates = @{
    "H42" = 50000
    "J42" = 50000
    "L42" = 50000
    "N42" = -50970
    "H115" = 50000
    "J115" = 50000
    "L115" = 50000
    "N115" = -52350
    "H122" = 14000
    "J122" = 8500
    "L122" = 25500
    "N122" = -30400
}
foreach ($cell in $updates.Keys) {
    $ws.Range($cell).Value = $updates[$cell]
}

$ws = $wb.Worksheets.Item("LTW")
$updates = @{
    "H22" = 10102452
    "J22" = 1830.6666
    "L22" = 1830.6666
    "N22" = -2420.6666
    "H27" = 10102452
    "J27" = 1830.6666
    "L27" = 1830.6666
    "N27" = -2044.6666
    "H55" = 912.6
    "I55" = 912.6
    "K55" = 912.6
    "M55" = -739.6
    "H93" = 3817.2144
    "I93" = 3929.1
    "K93" = 3929.1
    "M93" = -2681.1
}
foreach ($cell in $updates.Keys) {
    $ws.Range($cell).Value = $updates[$cell]
}

$ws = $wb.Worksheets.Item("WVR")
$updates = @{
    "H45" = 13909.223
    "I45" = 7968
    "K45" = 7968
    "M45" = -7477
    "H47" = 36438.59
    "I47" = 35464.066
    "J47" = 43747.5
    "K47" = 35464.066
    "L47" = 43747.5
    "M47" = -34892.066
    "N47" = -44891.5
    "H132" = 1652.125
    "I132" = 888.1429000000001
    "J132" = 7000
    "K132" = 2664.4287
    "L132" = 21000
    "M132" = -134.4287000000004
    "N132" = -26060
    "H136" = 1250
    "J136" = 2000
    "L136" = 6000
    "N136" = -11100
}
foreach ($cell in $updates.Keys) {
    $ws.Range($cell).Value = $updates[$cell]
}
